$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 6).Value = 3.35
$ws.Cells.Item(3, 7).Value = 4.1
$ws.Cells.Item(3, 9).Value = 2.42
$ws.Cells.Item(3, 15).Value = 1.29
$ws.Cells.Item(3, 16).Value = 1.81
$ws.Cells.Item(3, 17).Value = 2
$ws.Cells.Item(3, 19).Value = 3.6
$ws.Cells.Item(3, 20).Value = 1.78
$ws.Cells.Item(3, 21).Value = 2.02
$ws.Cells.Item(3, 22).Value = 1.7
$ws.Cells.Item(3, 23).Value = 1.33
$ws.Cells.Item(3, 38).Value = 70
$ws.Cells.Item(4, 8).Value = 4.7
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(4, 10).Value = 3.25
$ws.Cells.Item(4, 11).Value = 5.2
$ws.Cells.Item(4, 14).Value = 2.48
$ws.Cells.Item(4, 16).Value = 1.58
$ws.Cells.Item(4, 17).Value = 2.4
$ws.Cells.Item(4, 21).Value = 1.01
$ws.Cells.Item(5, 7).Value = 2.84
$ws.Cells.Item(5, 8).Value = 2.54
$ws.Cells.Item(5, 9).Value = 2.9
$ws.Cells.Item(5, 22).Value = 1.52
$ws.Cells.Item(5, 23).Value = 1.54
$ws.Cells.Item(6, 7).Value = 1.28
$ws.Cells.Item(6, 10).Value = 6.6
$ws.Cells.Item(6, 18).Value = 1.53
$ws.Cells.Item(6, 19).Value = 2.04
$ws.Cells.Item(6, 23).Value = 4.6
$ws.Cells.Item(7, 9).Value = 1.42
$ws.Cells.Item(7, 10).Value = 4.7
$ws.Cells.Item(7, 14).Value = 3.3
$ws.Cells.Item(7, 15).Value = 1.26
$ws.Cells.Item(7, 17).Value = 1.79
$ws.Cells.Item(7, 19).Value = 3.05
$ws.Cells.Item(7, 22).Value = 3.3
$ws.Cells.Item(8, 6).Value = 3.2
$ws.Cells.Item(8, 7).Value = 3.65
$ws.Cells.Item(8, 9).Value = 2.6
$ws.Cells.Item(8, 14).Value = 3.2
$ws.Cells.Item(8, 17).Value = 2.12
$ws.Cells.Item(8, 19).Value = 3.8
$ws.Cells.Item(8, 23).Value = 1.37
$ws.Cells.Item(8, 26).Value = 18
$ws.Cells.Item(8, 29).Value = 9
$ws.Cells.Item(8, 30).Value = 13.5
$ws.Cells.Item(8, 33).Value = 17
$ws.Cells.Item(8, 38).Value = 65
$ws.Cells.Item(9, 16).Value = 2.18
$ws.Cells.Item(9, 18).Value = 1.48
$ws.Cells.Item(9, 21).Value = 2.44
$ws.Cells.Item(11, 16).Value = 1.81
$ws.Cells.Item(11, 17).Value = 1.87
$ws.Cells.Item(11, 18).Value = 1.31
$ws.Cells.Item(11, 22).Value = 1.2
$ws.Cells.Item(12, 7).Value = 4.1
$ws.Cells.Item(12, 9).Value = 2.36
$ws.Cells.Item(12, 13).Value = 1.05
$ws.Cells.Item(12, 14).Value = 4
$ws.Cells.Item(12, 15).Value = 1.2
$ws.Cells.Item(12, 16).Value = 2.2
$ws.Cells.Item(12, 17).Value = 1.59
$ws.Cells.Item(12, 18).Value = 1.48
$ws.Cells.Item(12, 19).Value = 2.46
$ws.Cells.Item(12, 20).Value = 1.6
$ws.Cells.Item(12, 21).Value = 2.3
$ws.Cells.Item(12, 22).Value = 1.73
$ws.Cells.Item(12, 23).Value = 1.32
$ws.Cells.Item(12, 24).Value = 980
$ws.Cells.Item(12, 25).Value = 15
$ws.Cells.Item(12, 26).Value = 980
$ws.Cells.Item(12, 27).Value = 980
$ws.Cells.Item(12, 28).Value = 980
$ws.Cells.Item(12, 29).Value = 11
$ws.Cells.Item(12, 30).Value = 13.5
$ws.Cells.Item(12, 31).Value = 980
$ws.Cells.Item(12, 32).Value = 980
$ws.Cells.Item(12, 33).Value = 980
$ws.Cells.Item(12, 34).Value = 980
$ws.Cells.Item(12, 35).Value = 980
$ws.Cells.Item(12, 36).Value = 75
$ws.Cells.Item(12, 37).Value = 980
$ws.Cells.Item(12, 38).Value = 980
$ws.Cells.Item(12, 39).Value = 85
$ws.Cells.Item(12, 40).Value = 980
$ws.Cells.Item(12, 41).Value = 15.5
$ws.Cells.Item(13, 7).Value = 1.33
$ws.Cells.Item(14, 7).Value = 3.4
$ws.Cells.Item(14, 23).Value = 1.41
$ws.Cells.Item(15, 6).Value = 3.25
$ws.Cells.Item(15, 9).Value = 2.3
$ws.Cells.Item(15, 10).Value = 3.9
$ws.Cells.Item(15, 14).Value = 3.35
$ws.Cells.Item(15, 16).Value = 3.35
$ws.Cells.Item(15, 17).Value = 1.3
$ws.Cells.Item(15, 18).Value = 1.86
$ws.Cells.Item(15, 19).Value = 1.8
$ws.Cells.Item(15, 22).Value = 1.77
$ws.Cells.Item(17, 6).Value = 1.99
$ws.Cells.Item(17, 7).Value = 2.46
$ws.Cells.Item(17, 9).Value = 4.8
$ws.Cells.Item(17, 11).Value = 5.3
$ws.Cells.Item(17, 22).Value = 1.26
$ws.Cells.Item(17, 23).Value = 1.68
$ws.Cells.Item(18, 12).Value = 1.43
$ws.Cells.Item(18, 14).Value = 3.5
$ws.Cells.Item(18, 16).Value = 1.85
$ws.Cells.Item(18, 20).Value = 2.34
$ws.Cells.Item(18, 27).Value = 380
$ws.Cells.Item(18, 41).Value = 280
$ws.Cells.Item(19, 6).Value = 2.44
$ws.Cells.Item(19, 7).Value = 2.68
$ws.Cells.Item(19, 9).Value = 3.9
$ws.Cells.Item(19, 22).Value = 1.34
$ws.Cells.Item(19, 23).Value = 1.59
$ws.Cells.Item(20, 9).Value = 4.5
$ws.Cells.Item(20, 14).Value = 3.35
$ws.Cells.Item(20, 17).Value = 2.22
$ws.Cells.Item(20, 18).Value = 1.29
$ws.Cells.Item(20, 23).Value = 1.92
$ws.Cells.Item(21, 27).Value = 15
$ws.Cells.Item(21, 28).Value = 980
$ws.Cells.Item(22, 6).Value = 2.08
$ws.Cells.Item(22, 7).Value = 2.28
$ws.Cells.Item(22, 8).Value = 3.4
$ws.Cells.Item(22, 9).Value = 3.9
$ws.Cells.Item(22, 10).Value = 3.45
$ws.Cells.Item(22, 11).Value = 3.9
$ws.Cells.Item(22, 17).Value = 1.89
$ws.Cells.Item(22, 20).Value = 1.72
$ws.Cells.Item(22, 21).Value = 2.14
$ws.Cells.Item(22, 22).Value = 1.34
$ws.Cells.Item(22, 23).Value = 1.79
$ws.Cells.Item(22, 30).Value = 16
$ws.Cells.Item(22, 31).Value = 50
$ws.Cells.Item(22, 32).Value = 15
$ws.Cells.Item(22, 36).Value = 32
$ws.Cells.Item(22, 37).Value = 29
$ws.Cells.Item(22, 41).Value = 50
